$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the Price/Volume columns as text before writing, so values like
# "217.04" or "1.01" are preserved as literal strings (matching the source data)
# rather than being reinterpreted as numbers. ClearFormats() afterwards removes
# the temporary text-format styling again so cell styling is left untouched.
$fmtRange = $ws.Range("D2:E51")
$fmtRange.NumberFormat = "@"

$updates = @(
    @('D2', '26.702.04'),
    @('E2', '  -0.23%  '),
    @('D3', '1.635.68'),
    @('E3', '  -0.77%  '),
    @('E4', '  -0.03%  '),
    @('D5', '217.04'),
    @('E5', '  +0.33%  '),
    @('E6', '  -0.96%  '),
    @('E7', '  -0.04%  '),
    @('E8', '  -0.93%  '),
    @('E9', '  -0.84%  '),
    @('E10', '  -0.77%  '),
    @('D11', '0.0842'),
    @('E11', '  -0.04%  '),
    @('D12', '1.863.22'),
    @('E12', '  -0.83%  '),
    @('D13', '1.623.76'),
    @('E13', '  -1.75%  '),
    @('E14', '  -1.20%  '),
    @('D15', '0.524'),
    @('E15', '  -1.59%  '),
    @('D16', '64.35'),
    @('E16', '  -1.77%  '),
    @('D17', '26.696.95'),
    @('E17', '  -0.35%  '),
    @('E18', '  -2.79%  '),
    @('B19', 'BitcoinCash'),
    @('C19', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'),
    @('D19', '211.02'),
    @('E19', '  -3.58%  '),
    @('B20', 'Dai'),
    @('C20', 'https://coinranking.com/coin/MoTuySvg7+dai-dai'),
    @('D20', '1.01'),
    @('E20', '  -0.03%  '),
    @('E21', '  -1.03%  '),
    @('D22', '6.18'),
    @('E22', '  -1.46%  '),
    @('D23', '2.29'),
    @('E23', '  -3.03%  '),
    @('D24', '9.25'),
    @('E24', '  -2.83%  '),
    @('D25', '145.93'),
    @('E25', '  -0.33%  '),
    @('E26', '  -0.12%  '),
    @('D27', '0.118'),
    @('E27', '  -1.96%  '),
    @('D28', '7.06'),
    @('E28', '  -0.94%  '),
    @('D29', '15.53'),
    @('E29', '  -1.38%  '),
    @('E30', '  -2.61%  '),
    @('E31', '  +0.54%  '),
    @('E32', '  -0.27%  '),
    @('E33', '  -1.53%  '),
    @('D34', '1.271.02'),
    @('E34', '  -0.81%  '),
    @('E35', '  -1.47%  '),
    @('E36', '  +0.12%  '),
    @('E37', '  -2.13%  '),
    @('D38', '0.527'),
    @('E38', '  -1.78%  '),
    @('D39', '0.806'),
    @('E39', '  -2.75%  '),
    @('E40', '  -0.05%  '),
    @('D41', '0.801'),
    @('E41', '  -1.62%  '),
    @('E42', '  -2.54%  '),
    @('B43', 'FraxShare'),
    @('C43', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'),
    @('D43', '5.26'),
    @('E43', '  -3.65%  '),
    @('B44', 'RocketPoolETH'),
    @('C44', 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'),
    @('D44', '1.773.94'),
    @('E44', '  -0.89%  '),
    @('D45', '91.22'),
    @('E45', '  -0.86%  '),
    @('D46', '60.26'),
    @('E46', '  +0.76%  '),
    @('E47', '  -1.97%  '),
    @('B48', 'Cronos'),
    @('C48', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'),
    @('D48', '0.0518'),
    @('E48', '  +0.36%  '),
    @('B49', 'EnergySwap'),
    @('C49', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'),
    @('D49', '7.52'),
    @('E49', '  -3.45%  '),
    @('D50', '0.0960'),
    @('E50', '  -1.11%  '),
    @('B51', 'USDD'),
    @('C51', 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'),
    @('D51', '1.01'),
    @('E51', '  -0.06%  ')
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

$fmtRange.ClearFormats()

Write-Output "Applied $($updates.Count) cell updates"
